$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.028799999999987
$ws.Range("B14").Value = 9.124100000000004
$ws.Range("B16").Value = 9.495200000000004
$ws.Range("B21").Value = 5.87739999999999
$ws.Range("B23").Value = 5.214800000000001
$ws.Range("B25").Value = 5.89929999999999
